$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.1443736666666667
$ws.Range("H2").Value = 0.433121
$ws.Range("I2").Value = 0.7378778224885942
$ws.Range("J2").Value = 0.7378778224885942
$ws.Range("M2").Value = 2.113523666666667
$ws.Range("N2").Value = 6.340571000000001
$ws.Range("O2").Value = 0.2651220308693004
$ws.Range("P2").Value = 0.2651220308693004
$ws.Range("Q2").Value = 0.3051371613434445
$ws.Range("R2").Value = 2.746234452091
$ws.Range("S2").Value = 0.1956276668315932
$ws.Range("T2").Value = 0.1956276668315932

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.1443736666666667
$ws.Range("H3").Value = 0.433121
$ws.Range("I3").Value = 0.7378778224885942
$ws.Range("J3").Value = 0.7378778224885942
$ws.Range("O3").Value = 0.2869289465860668
$ws.Range("P3").Value = 0.2869289465860668
$ws.Range("Q3").Value = 0.3302354164286667
$ws.Range("R3").Value = 2.972118747858
$ws.Range("S3").Value = 0.2117185063158731
$ws.Range("T3").Value = 0.2117185063158731

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.1443736666666667
$ws.Range("H4").Value = 0.433121
$ws.Range("I4").Value = 0.7378778224885942
$ws.Range("J4").Value = 0.7378778224885942
$ws.Range("M4").Value = 1.164746666666667
$ws.Range("N4").Value = 3.49424
$ws.Range("O4").Value = 0.1461067158059967
$ws.Range("P4").Value = 0.1461067158059966
$ws.Range("Q4").Value = 0.1681587470044444
$ws.Range("R4").Value = 1.51342872304
$ws.Range("S4").Value = 0.1078089053098887
$ws.Range("T4").Value = 0.1078089053098887

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.1443736666666667
$ws.Range("H5").Value = 0.433121
$ws.Range("I5").Value = 0.7378778224885942
$ws.Range("J5").Value = 0.7378778224885942
$ws.Range("M5").Value = 2.406253666666667
$ws.Range("N5").Value = 7.218761
$ws.Range("O5").Value = 0.3018423067386362
$ws.Range("P5").Value = 0.3018423067386362
$ws.Range("Q5").Value = 0.3473996647867778
$ws.Range("R5").Value = 3.126596983081
$ws.Range("S5").Value = 0.2227227440312392
$ws.Range("T5").Value = 0.2227227440312392

# Row 6
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.3333333333333333
$ws.Range("G6").Value = 0.051287
$ws.Range("H6").Value = 0.153861
$ws.Range("I6").Value = 0.2621221775114058
$ws.Range("J6").Value = 0.2621221775114058
$ws.Range("M6").Value = 2.113523666666667
$ws.Range("N6").Value = 6.340571000000001
$ws.Range("O6").Value = 0.2651220308693004
$ws.Range("P6").Value = 0.2651220308693004
$ws.Range("Q6").Value = 0.1083962882923333
$ws.Range("R6").Value = 0.975566594631
$ws.Range("S6").Value = 0.06949436403770715
$ws.Range("T6").Value = 0.06949436403770716

# Row 7
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0.3333333333333333
$ws.Range("G7").Value = 0.051287
$ws.Range("H7").Value = 0.153861
$ws.Range("I7").Value = 0.2621221775114058
$ws.Range("J7").Value = 0.2621221775114058
$ws.Range("O7").Value = 0.2869289465860668
$ws.Range("P7").Value = 0.2869289465860668
$ws.Range("Q7").Value = 0.117312140042
$ws.Range("R7").Value = 1.055809260378
$ws.Range("S7").Value = 0.07521044027019366
$ws.Range("T7").Value = 0.07521044027019368

# Row 8
$ws.Range("E8").Value = 1
$ws.Range("F8").Value = 0.3333333333333333
$ws.Range("G8").Value = 0.051287
$ws.Range("H8").Value = 0.153861
$ws.Range("I8").Value = 0.2621221775114058
$ws.Range("J8").Value = 0.2621221775114058
$ws.Range("M8").Value = 1.164746666666667
$ws.Range("N8").Value = 3.49424
$ws.Range("O8").Value = 0.1461067158059967
$ws.Range("P8").Value = 0.1461067158059966
$ws.Range("Q8").Value = 0.05973636229333333
$ws.Range("R8").Value = 0.53762726064
$ws.Range("S8").Value = 0.03829781049610797
$ws.Range("T8").Value = 0.03829781049610797

# Row 9
$ws.Range("E9").Value = 1
$ws.Range("F9").Value = 0.3333333333333333
$ws.Range("G9").Value = 0.051287
$ws.Range("H9").Value = 0.153861
$ws.Range("I9").Value = 0.2621221775114058
$ws.Range("J9").Value = 0.2621221775114058
$ws.Range("M9").Value = 2.406253666666667
$ws.Range("N9").Value = 7.218761
$ws.Range("O9").Value = 0.3018423067386362
$ws.Range("P9").Value = 0.3018423067386362
$ws.Range("Q9").Value = 0.1234095318023333
$ws.Range("R9").Value = 1.110685786221
$ws.Range("S9").Value = 0.07911956270739699
$ws.Range("T9").Value = 0.07911956270739699
